$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-08-31 07:14:38"
$newTimestamp = "2022-08-31 21:01:28"

for ($row = 2; $row -le 65; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # Column O is the 15th column
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
